$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.95
$ws.Range("I2").Value = 4.1
$ws.Range("J2").Value = 2.75
$ws.Range("L2").Value = 4.75
$ws.Range("Z2").Value = 17
$ws.Range("AL2").Value = 41
$ws.Range("AZ2").Value = 81
